# repull data, push all data, mean calculation
# Update the dSF column (F) values to reflect repulled/recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 2
    8  = -3
    12 = 3
    19 = -2
    20 = 7
    23 = 0
    25 = 1
    26 = 1
    31 = 3
    33 = 0
    35 = 0
    38 = -2
    40 = 1
    47 = -1
    50 = 2
    54 = 0
    62 = 2
    68 = 4
    76 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
